$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header (C1) with the new date (copy B1 formatting so it
# reuses the same bold/bordered/centered style)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# Reorder rows 2-8: move the 5 named rows up (A4:B8 -> A2:B6),
# then place avg/total at the bottom (rows 7-8), and fill column C values.

$ws.Range("A2").Value = "Alpha Mega"
$ws.Range("B2").Value = 57248.24
$ws.Range("C2").Value = 58970.63

$ws.Range("A3").Value = "Alpha Mercosur"
$ws.Range("B3").Value = 5262.3
$ws.Range("C3").Value = 5448.26

$ws.Range("A4").Value = "Arpenta acciones"
$ws.Range("B4").Value = 3495.3
$ws.Range("C4").Value = 3506.39

$ws.Range("A5").Value = "Fima Acciones"
$ws.Range("B5").Value = 27544.64
$ws.Range("C5").Value = 24413.85

$ws.Range("A6").Value = "Fima PB Acciones"
$ws.Range("B6").Value = 41469.66
$ws.Range("C6").Value = 41245.14

$ws.Range("A7").Value = "avg"
$ws.Range("B7").Value = 27004.03
$ws.Range("C7").Value = 26716.85

$ws.Range("A8").Value = "total"
$ws.Range("B8").Value = 135020.14
$ws.Range("C8").Value = 133584.27
